$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.009821333333333333
$ws.Range("H2").Value = 0.029464
$ws.Range("I2").Value = 0.06297798848338983
$ws.Range("J2").Value = 0.06297798848338984
$ws.Range("M2").Value = 0.240998
$ws.Range("N2").Value = 0.7229939999999999
$ws.Range("O2").Value = 0.05495977716704094
$ws.Range("P2").Value = 0.05495977716704094
$ws.Range("Q2").Value = 0.002366921690666666
$ws.Range("R2").Value = 0.021302295216
$ws.Range("S2").Value = 0.003461256213475576
$ws.Range("T2").Value = 0.003461256213475576

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.009821333333333333
$ws.Range("H3").Value = 0.029464
$ws.Range("I3").Value = 0.06297798848338983
$ws.Range("J3").Value = 0.06297798848338984
$ws.Range("O3").Value = 0.3551990176181375
$ws.Range("P3").Value = 0.3551990176181375
$ws.Range("Q3").Value = 0.01529715553155555
$ws.Range("R3").Value = 0.137674399784
$ws.Range("S3").Value = 0.02236971964086644
$ws.Range("T3").Value = 0.02236971964086645

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.009821333333333333
$ws.Range("H4").Value = 0.029464
$ws.Range("I4").Value = 0.06297798848338983
$ws.Range("J4").Value = 0.06297798848338984
$ws.Range("M4").Value = 2.586447
$ws.Range("N4").Value = 7.759341
$ws.Range("O4").Value = 0.5898412052148215
$ws.Range("P4").Value = 0.5898412052148215
$ws.Range("Q4").Value = 0.025402358136
$ws.Range("R4").Value = 0.228621223224
$ws.Range("S4").Value = 0.03714701262904781
$ws.Range("T4").Value = 0.03714701262904781

# Row 5
$ws.Range("I5").Value = 0.3247949111459754
$ws.Range("J5").Value = 0.3247949111459754
$ws.Range("M5").Value = 0.240998
$ws.Range("N5").Value = 0.7229939999999999
$ws.Range("O5").Value = 0.05495977716704094
$ws.Range("P5").Value = 0.05495977716704094
$ws.Range("Q5").Value = 0.01220687003066667
$ws.Range("R5").Value = 0.109861830276
$ws.Range("S5").Value = 0.01785065594157167
$ws.Range("T5").Value = 0.01785065594157167

# Row 6
$ws.Range("I6").Value = 0.3247949111459754
$ws.Range("J6").Value = 0.3247949111459754
$ws.Range("O6").Value = 0.3551990176181375
$ws.Range("P6").Value = 0.3551990176181375
$ws.Range("S6").Value = 0.1153668333664207
$ws.Range("T6").Value = 0.1153668333664207

# Row 7
$ws.Range("I7").Value = 0.3247949111459754
$ws.Range("J7").Value = 0.3247949111459754
$ws.Range("M7").Value = 2.586447
$ws.Range("N7").Value = 7.759341
$ws.Range("O7").Value = 0.5898412052148215
$ws.Range("P7").Value = 0.5898412052148215
$ws.Range("Q7").Value = 0.131006989146
$ws.Range("R7").Value = 1.179062902314
$ws.Range("S7").Value = 0.191577421837983
$ws.Range("T7").Value = 0.191577421837983

# Row 8
$ws.Range("G8").Value = 0.09547600000000001
$ws.Range("H8").Value = 0.286428
$ws.Range("I8").Value = 0.6122271003706348
$ws.Range("J8").Value = 0.6122271003706349
$ws.Range("M8").Value = 0.240998
$ws.Range("N8").Value = 0.7229939999999999
$ws.Range("O8").Value = 0.05495977716704094
$ws.Range("P8").Value = 0.05495977716704094
$ws.Range("Q8").Value = 0.023009525048
$ws.Range("R8").Value = 0.207085725432
$ws.Range("S8").Value = 0.0336478650119937
$ws.Range("T8").Value = 0.0336478650119937

# Row 9
$ws.Range("G9").Value = 0.09547600000000001
$ws.Range("H9").Value = 0.286428
$ws.Range("I9").Value = 0.6122271003706348
$ws.Range("J9").Value = 0.6122271003706349
$ws.Range("O9").Value = 0.3551990176181375
$ws.Range("P9").Value = 0.3551990176181375
$ws.Range("Q9").Value = 0.1487080391186667
$ws.Range("R9").Value = 1.338372352068
$ws.Range("S9").Value = 0.2174624646108503
$ws.Range("T9").Value = 0.2174624646108504

# Row 10
$ws.Range("G10").Value = 0.09547600000000001
$ws.Range("H10").Value = 0.286428
$ws.Range("I10").Value = 0.6122271003706348
$ws.Range("J10").Value = 0.6122271003706349
$ws.Range("M10").Value = 2.586447
$ws.Range("N10").Value = 7.759341
$ws.Range("O10").Value = 0.5898412052148215
$ws.Range("P10").Value = 0.5898412052148215
$ws.Range("Q10").Value = 0.246943613772
$ws.Range("R10").Value = 2.222492523948
$ws.Range("S10").Value = 0.3611167707477907
$ws.Range("T10").Value = 0.3611167707477908
